$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update parameter search ranges
$ws.Range("C2").Value = "UNIF(1000 - 10000)"
$ws.Range("C5").Value = "10^UNIF(-8 - -7)"
$ws.Range("C6").Value = "10^UNIF(-10 - -6) or UNIF(0 - 0.25)/(4*N)"
$ws.Range("C10").Value = "if growing: UNIF(1.01 - 2)*N ; if shrinking: 1/UNIF(1.01-2)*N ; if 2-cycling: N; if chaotic: N"
$ws.Range("C12").Value = "0.5: 1; 0.5: UNIF(0.75 - 1)"
$ws.Range("C16").Value = "10^UNIF(1 - 4)"

# Fix time tracking: move the active selection from C25 to C24
$ws.Range("C24").Select()
